$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column (C) for rows 2-5
# from serial 45208 (2023-10-09) to serial 45212 (2023-10-13)
foreach ($r in 2..5) {
    $ws.Cells.Item($r, 3).Value = 45212
}
